$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Admin"
$ws.Range("B5").Value = "admin123"

$ws.Range("B5").Select()
